$d = $word.ActiveDocument

# Locate the target sentence (the run containing it currently ends with
# "...making an impact in IT Sector. " - trailing period + space).
$target = $d.Content.Duplicate
$found = $target.Find.Execute("execute my vision of making an impact in IT Sector. ")
if (-not $found) {
    throw "Could not find the target sentence to edit."
}

$sentenceStart = $target.Start
$sentenceEnd = $target.End
$sentenceText = $target.Text

# Work out where to split the run: right before " making an impact in IT
# Sector" (the tail of the original run text), and where the trailing
# ". " needs to be dropped.
$splitMarker = " making an impact in IT Sector"
$markerIdx = $sentenceText.IndexOf($splitMarker)
if ($markerIdx -lt 0) {
    throw "Could not locate the split point inside the target sentence."
}

$splitPos = $sentenceStart + $markerIdx
$trailing = ". "
$tailStart = $sentenceEnd - $trailing.Length

# Sanity-check that the text we are about to delete really is the
# trailing ". " we expect.
$trailingRange = $d.Range($tailStart, $sentenceEnd)
if ($trailingRange.Text -ne $trailing) {
    throw "Unexpected trailing text: [$($trailingRange.Text)]"
}

# Remove the trailing period + space - the sentence now ends with
# "...making an impact in IT Sector" (no period, no trailing space).
$trailingRange.Delete()
$newEnd = $tailStart

# Split the single run into two runs at $splitPos..$newEnd by re-pasting
# its own formatted text onto itself - this preserves the run's existing
# formatting (rPr) exactly while forcing a run boundary at $splitPos.
$tailRange = $d.Range($splitPos, $newEnd)
$tailCopy = $tailRange.Duplicate()
$tailRange.FormattedText = $tailCopy.FormattedText
